$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sessions")

# Insert a new row at row 19 (pushes existing row 19 "Context New Vehicle
# Market File" and everything below it down by one row).
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the "Context Fuel Upstream File" parameter,
# matching the formatting/style already used by the surrounding rows
# (the Insert() above carries the row-18 formatting down automatically).
$ws.Range("A19").Value = "Context Fuel Upstream File"
$ws.Range("B19").Value = "String"
$ws.Range("C19").Value = "input_samples/context_fuel_upstream.csv"
$ws.Range("D19").Value = "input_samples/context_fuel_upstream.csv"

$ws.Range("C19").Select()
